# BOM.xlsx update: added JLPCB gerber files
# Swaps a handful of resistor / LED BOM lines over to new manufacturers &
# part numbers (Vishay -> Stackpole on the 10k, Lite-On -> Wurth on 3 of
# the 4 LEDs, designator/qty reshuffle on the R3/R4/R5/R6 resistor bank).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 20 (Item 14, R2): 10k resistor, new manufacturer/part/desc ---
$ws.Range("D20").Value = "Stackpole Electronics Inc"
$ws.Range("E20").Value = "RNCP0805FTD10K0"
$ws.Range("F20").Value = "RES 10K OHM 1% 1/4W 0805"

# --- Row 21 (Item 15): was "R3, R5" qty 2 -> now just "R5" qty 1 ---
$ws.Range("B21").Value = "R5"
$ws.Range("C21").Value = 1

# --- Row 22 (Item 16): was "R4" -> now "R3, R4"; resistor value 4.7R -> 5.6R 1% ---
$ws.Range("B22").Value = "R3, R4"
$ws.Range("E22").Value = "CRCW08055R60FKEA"
$ws.Range("F22").Value = "RES SMD 5.6 OHM 1% 1/8W 0805"

# --- Row 23 (Item 17, R6): new part number + wattage in description ---
$ws.Range("E23").Value = "RCS0805560RJNEA"
$ws.Range("F23").Value = "RES SMD 560 OHM 5% 0.4W 0805"

# --- Row 34 (Item 28, L1 - red LED): Lite-On -> Wurth Elektronik ---
$ws.Range("D34").Value = "Würth Elektronik"
$ws.Range("E34").Value = "150080RS75000"
$ws.Range("F34").Value = "LED RED CLEAR 0805 SMD"

# --- Row 35 (Item 29, L2 - green LED): Lite-On -> Wurth Elektronik ---
$ws.Range("D35").Value = "Würth Elektronik"
$ws.Range("E35").Value = "150080GS75000"
$ws.Range("F35").Value = "LED GREEN CLEAR 0805 SMD"

# --- Row 36 (Item 30, L3 - blue LED): Lite-On -> Wurth Elektronik ---
$ws.Range("D36").Value = "Würth Elektronik"
$ws.Range("E36").Value = "150080BS75000"
$ws.Range("F36").Value = "LED BLUE CLEAR 0805 SMD"

# --- Row 37 (Item 31, L4 - orange LED): stays Lite-On, new part + package ---
$ws.Range("E37").Value = "LTST-C170KFKT"
$ws.Range("G37").Value = "0805"

# --- Restore the view state (scroll position / selection) from the edit ---
$ws.Range("F23:F24").Select()
$excel.ActiveWindow.ScrollRow = 11
